$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 13, shrinking the used range down to A1:C3
$ws.Range("A4:C13").EntireRow.Delete()

# Update row 2 values
$ws.Range("A2").Value = "5001327-37.2021.8.21.0085"
$ws.Range("B2").Value = "9000091-79.2021.8.21.0085"
$ws.Range("C2").Value = "Migrado"

# Update row 3 values
$ws.Range("A3").Value = "5003574-61.2017.8.21.0010"
$ws.Range("B3").Value = "0020708-89.2017.8.21.0010"
$ws.Range("C3").Value = "Digitalizado"
